$d = $word.ActiveDocument

# 1) Replace the signature name: "Dinesh S" -> "Dibyaranjan Dalai"
#    (the two runs "Dinesh" + " S" collapse into one run with the new name)
$d.Content.Find.Execute("Dinesh S", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Dibyaranjan Dalai", 2)

# 2) Replace the title: "Senior Consultant" -> "Principal Consultant"
$d.Content.Find.Execute("Senior Consultant", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Principal Consultant", 2)

# 3) Split "Phone:" into two runs "Phone" + ":" by re-typing the colon
#    with its own formatting (keeps same visible text, different run split).
$rng = $d.Content
$rng.Find.Execute("Phone:", $true, $false, $false, $false, $false,
                   $true, 1, $false, "", 0)
if ($rng.Find.Found) {
    $colonStart = $rng.End - 1
    $colonRange = $d.Range($colonStart, $rng.End)
    $colonRange.Text = ":"
}

# 4) Replace the phone number
$d.Content.Find.Execute("+91-8618520409", $true, $false, $false, $false, $false,
                         $true, 1, $false, "+1-7034594554", 2)
